# "solved problem about gender setting"
# The Customer sheet had a bad row (id=5, yolanda, gender FEMALE, ...) at row 6.
# That row is removed (shifting the row below it up), and a new customer row
# (id=7, yyy, gender "null", isVip 0, loyaltyCard 0, password 123456,
# phone 12234567890) is appended at the new last row (row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer")

# Remove the old row 6 entirely; row 7 shifts up to become the new row 6.
$ws.Rows(6).Delete()

# Append the brand-new customer record as the new row 7.
$ws.Cells.Item(7, 1).Value = "7"
$ws.Cells.Item(7, 2).Value = "yyy"
$ws.Cells.Item(7, 3).Value = "null"
$ws.Cells.Item(7, 4).Value = "0"
$ws.Cells.Item(7, 5).Value = "0"
$ws.Cells.Item(7, 6).Value = "123456"
$ws.Cells.Item(7, 7).Value = "12234567890"
